$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 121341.1328125
$ws.Range("H4").Value = 20.93000030517578
$ws.Range("H5").Value = 11.27999973297119
$ws.Range("H6").Value = 47.40999984741211
$ws.Range("H7").Value = 30.04000091552734
$ws.Range("H8").Value = 49.59000015258789
$ws.Range("B9").Value = 10.85123538970947
$ws.Range("C9").Value = 11.32734107971191
$ws.Range("D9").Value = 11.0496129989624
$ws.Range("E9").Value = 10.61318302154541
$ws.Range("F9").Value = 10.66277694702148
$ws.Range("G9").Value = 10.28586006164551
$ws.Range("H9").Value = 10.39000034332275
$ws.Range("H10").Value = 88.98999786376953
$ws.Range("C11").Value = 27.02092742919922
$ws.Range("H11").Value = 26.61000061035156
$ws.Range("B12").Value = 13.35010242462158
$ws.Range("C12").Value = 13.41001319885254
$ws.Range("D12").Value = 12.95069789886475
$ws.Range("E12").Value = 12.65114402770996
$ws.Range("F12").Value = 12.94999980926514
$ws.Range("G12").Value = 12.82999992370605
$ws.Range("H12").Value = 12.39999961853027
$ws.Range("H13").Value = 24.93000030517578
$ws.Range("H14").Value = 32.47999954223633
$ws.Range("H15").Value = 20.54999923706055
$ws.Range("H16").Value = 12.02999973297119
$ws.Range("H17").Value = 1.740000009536743
$ws.Range("H18").Value = 12.89000034332275
$ws.Range("H19").Value = 1.899999976158142
$ws.Range("H20").Value = 14.78999996185303
$ws.Range("H21").Value = 37.86999893188477
$ws.Range("H22").Value = 12.89999961853027
$ws.Range("H23").Value = 14.21000003814697
$ws.Range("H24").Value = 17.73999977111816
$ws.Range("H25").Value = 7.150000095367432
$ws.Range("H26").Value = 3.769999980926514
$ws.Range("H27").Value = 110.0800018310547
$ws.Range("H28").Value = 19.71999931335449
$ws.Range("H29").Value = 98.36000061035156
$ws.Range("H30").Value = 82.62000274658203
$ws.Range("B31").Value = 9.919038772583008
$ws.Range("C31").Value = 10.13638782501221
$ws.Range("D31").Value = 9.84000301361084
$ws.Range("E31").Value = 9.612772941589355
$ws.Range("F31").Value = 9.606365203857422
$ws.Range("G31").Value = 9.477620124816895
$ws.Range("H31").Value = 9.680000305175781
$ws.Range("B32").Value = 32.374755859375
$ws.Range("C32").Value = 32.94986343383789
$ws.Range("D32").Value = 31.58149909973145
$ws.Range("E32").Value = 30.76841163635254
$ws.Range("F32").Value = 31.2026195526123
$ws.Range("G32").Value = 30.89505767822266
$ws.Range("H33").Value = 30.67000007629395
$ws.Range("H34").Value = 4.190000057220459
$ws.Range("H35").Value = 12.52999973297119
$ws.Range("H36").Value = 3.990000009536743
$ws.Range("H37").Value = 72.23999786376953
$ws.Range("H38").Value = 5.940000057220459
$ws.Range("H39").Value = 11.25
$ws.Range("H40").Value = 5.210000038146973
$ws.Range("H41").Value = 36.72000122070312
$ws.Range("H42").Value = 3.450000047683716
$ws.Range("H43").Value = 14.63000011444092
$ws.Range("H44").Value = 41.33000183105469
$ws.Range("H45").Value = 20.04999923706055
$ws.Range("H46").Value = 26.35000038146973
$ws.Range("H47").Value = 11.09000015258789
$ws.Range("H48").Value = 25.61000061035156
$ws.Range("H49").Value = 17.45999908447266
$ws.Range("B50").Value = 27.90649223327637
$ws.Range("C50").Value = 26.66357231140137
$ws.Range("D50").Value = 26.54519844055176
$ws.Range("E50").Value = 26.77207946777344
$ws.Range("F50").Value = 27.86703300476074
$ws.Range("G50").Value = 30.32327651977539
$ws.Range("H50").Value = 32.36999893188477
$ws.Range("H51").Value = 5.860000133514404
$ws.Range("H52").Value = 48.84999847412109
$ws.Range("H53").Value = 33.59999847412109
$ws.Range("H54").Value = 15.77999973297119
$ws.Range("H55").Value = 60.83000183105469
$ws.Range("H56").Value = 5.820000171661377
$ws.Range("H57").Value = 44.4900016784668
$ws.Range("H58").Value = 41.08000183105469
$ws.Range("H59").Value = 10.84000015258789
$ws.Range("H60").Value = 544.510009765625
$ws.Range("H61").Value = 89.63999938964844
$ws.Range("H62").Value = 41.33000183105469
$ws.Range("H63").Value = 122.9700012207031
$ws.Range("H64").Value = 228.4100036621094
$ws.Range("H65").Value = 77.66999816894531
$ws.Range("H66").Value = 69.38999938964844
$ws.Range("H67").Value = 146.8099975585938
$ws.Range("H68").Value = 182.6100006103516
$ws.Range("H69").Value = 77.04360198974609
$ws.Range("H70").Value = 66.01000213623047
$ws.Range("H71").Value = 30.18000030517578
$ws.Range("H72").Value = 42.66999816894531
$ws.Range("H73").Value = 93.95999908447266
$ws.Range("H74").Value = 214.7799987792969
